# The workbook originally holds two sheets: "login" (Username/Password test
# data) and an empty "user_details" sheet. This commit repurposes the file
# for a checkout/purchase-product test: the empty sheet is dropped, and the
# remaining sheet is renamed and repopulated with FirstName/LastName/
# PostalCode test data for one customer.

$wb = $excel.ActiveWorkbook

# Drop the unused second sheet.
$wb.Worksheets.Item("user_details").Delete()

# Reuse the remaining sheet, renaming it to match its new purpose.
$ws = $wb.Worksheets.Item("login")
$ws.Name = "CheckoutDetails"

# The old table was A1:B6 (2 cols x 6 rows); the new one is A1:C4
# (3 cols x 4 rows). Extend the existing bordered formatting into the new
# column C by copying it across from already-bordered cells, instead of
# building the border from scratch.
$ws.Range("A1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("A4").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("C4").PasteSpecial(-4122)

# Drop the now-unneeded trailing rows (old table had 6 rows, new one 4).
$ws.Rows("5:6").Delete()

# Recolor the header row's fill to the new solid green.
$ws.Range("A1:C1").Interior.Color = 5287936

# New header row.
$ws.Range("A1").Value = "FirstName"
$ws.Range("B1").Value = "LastName"
$ws.Range("C1").Value = "PostalCode"

# New data row; rows 3-4 stay blank (but keep their border/fill styling).
$ws.Range("A2").Value = "Selina"
$ws.Range("B2").Value = "Mayinga"
$ws.Range("C2").Value = 6500
$ws.Range("A3:C4").ClearContents()

# Widen the columns to fit the new (longer) headers.
$ws.Columns("A").ColumnWidth = 14.8
$ws.Columns("B").ColumnWidth = 16.4
$ws.Columns("C").ColumnWidth = 18.5

# Leave the selection on row 3, matching the saved UI state.
$ws.Range("A3:XFD3").Select()
